# Canpotex Data feed updated for 2023 Manual forecast
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated BrazilCFR forecast values (column B) ---
$ws.Range("B171").Value = 781.25
$ws.Range("B173").Value = 1173.75
$ws.Range("B176").Value = 968.75
$ws.Range("B178").Value = 731
$ws.Range("B179").Value = 629
$ws.Range("B180").Value = 564
$ws.Range("B181").Value = 516

# --- Restore the view/selection state left by the editor ---
try {
    $excel.ActiveWindow.ScrollRow = 165
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}

$ws.Range("E184").Select()
